# Fix uncertainty data in LCIs
# - Clear the now-unused "M" column placeholder zeros on rows 14-42
# - Row 43: H43 becomes a computed LN() formula instead of a hard-coded
#   literal, and M43 becomes a boolean TRUE flag instead of 0
# - Update the sheet view's selection/scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biogas")

# Remove the stray "0" placeholders in column M for rows 14 through 42
for ($row = 14; $row -le 42; $row++) {
    $ws.Cells.Item($row, 13).ClearContents()
}

# Row 43: H43 now computes its value from B43 via a formula
$ws.Range("H43").Formula = "=LN(B43*-1)"

# Row 43: M43 now holds a boolean TRUE instead of the numeric 0
$ws.Range("M43").Value = $true

# Reflect the author's on-screen selection/scroll position at save time
$ws.Range("H44").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
